$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 entirely (the original second data row), shifting rows 6:11 up to 5:10
$ws.Rows("5:5").Delete()

# Update selection to match the post-edit state (row 5 selected, whole row)
$ws.Range("A5:XFD5").Select()
